$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sensData")

# --- Row 2: "a" -> "Te" ---
$ws.Range("A2").Value = "Te"
# D2 loses its "=$J$1" formula and becomes a plain static value (0.3 unchanged)
$ws.Range("D2").Value = 0.3
# F2 drives B2/C2 (existing formulas F2*(1-D2) / F2*(1+D2)) to recalc to 350 / 650
$ws.Range("F2").Value = 500

# --- Row 3: "x" -> "pe" ---
$ws.Range("A3").Value = "pe"
$ws.Range("D3").Value = 0.3
$ws.Range("F3").Value = 300

# --- Row 4: new data row "pa" (create new shared strings in diff order: pa, group2, group3) ---
$ws.Range("A4").Value = "pa"
$ws.Range("B4").Formula = "=F4*(1-D4)"
$ws.Range("C4").Formula = "=F4*(1+D4)"
$ws.Range("D4").Value = 0.3
$ws.Range("F4").Value = 1000

# --- group labels: E3 group1 -> group2, E4 -> group3 (after A4 so "pa" gets index 23) ---
$ws.Range("E3").Value = "group2"
$ws.Range("E4").Value = "group3"

# --- Row 5: F5 becomes a present-but-empty cell (style reset to the default "Standard") ---
$ws.Range("F5").Style = "Standard"

$wb.Application.Calculate()

# Move the active selection to D3, matching the saved view state
$ws.Range("D3").Select()
